# Update stock report figures: swap the two duplicate-SKU rows' stats
# (qty/rate/variance/value) for several product pairs, and refresh the
# corresponding Sub Total / Grand Total rows to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B127").Value = 64329
$ws.Range("E127").Value = 128.32
$ws.Range("F127").Value = 2
$ws.Range("G127").Value = 241.38
$ws.Range("B128").Value = 57552
$ws.Range("E128").Value = 136.86
$ws.Range("F128").Value = -5
$ws.Range("G128").Value = -603.45
$ws.Range("F205").Value = 23
$ws.Range("G205").Value = 8674.219999999999
$ws.Range("B216").Value = 42397.92
$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6
$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 66
$ws.Range("G227").Value = 9522.48
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("B229").Value = 57802
$ws.Range("E229").Value = 162.71
$ws.Range("F229").Value = -79
$ws.Range("G229").Value = -11334.92
$ws.Range("B230").Value = 63531
$ws.Range("E230").Value = 152.53
$ws.Range("F230").Value = 63
$ws.Range("G230").Value = 9039.24
$ws.Range("B232").Value = 55356
$ws.Range("E232").Value = 54.04
$ws.Range("F232").Value = -158
$ws.Range("G232").Value = -7527.12
$ws.Range("B233").Value = 63510
$ws.Range("E233").Value = 50.66
$ws.Range("F233").Value = 117
$ws.Range("G233").Value = 5573.88
$ws.Range("F270").Value = 20
$ws.Range("G270").Value = 644.8
$ws.Range("B275").Value = 5353.95
$ws.Range("B366").Value = 53263
$ws.Range("E366").Value = 15.29
$ws.Range("F366").Value = -309
$ws.Range("G366").Value = -3958.29
$ws.Range("B367").Value = 65066
$ws.Range("E367").Value = 13.61
$ws.Range("F367").Value = 90
$ws.Range("G367").Value = 1152.9
$ws.Range("B375").Value = 64927
$ws.Range("E375").Value = 17.26
$ws.Range("F375").Value = 106
$ws.Range("G375").Value = 1719.32
$ws.Range("B376").Value = 45718
$ws.Range("E376").Value = 19.38
$ws.Range("F376").Value = -294
$ws.Range("G376").Value = -4768.68
$ws.Range("B380").Value = 45709
$ws.Range("E380").Value = 15.69
$ws.Range("F380").Value = -300
$ws.Range("G380").Value = -3945
$ws.Range("B381").Value = 64925
$ws.Range("E381").Value = 13.97
$ws.Range("F381").Value = 111
$ws.Range("G381").Value = 1459.65
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("B463").Value = 60025
$ws.Range("E463").Value = 37.22
$ws.Range("F463").Value = -98
$ws.Range("G463").Value = -3217.34
$ws.Range("B464").Value = 64833
$ws.Range("E464").Value = 34.9
$ws.Range("F464").Value = 95
$ws.Range("G464").Value = 3118.85
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 108
$ws.Range("G474").Value = 3545.64
$ws.Range("B619").Value = 1797704.03
$ws.Range("B620").Value = 1797704.03
